$d = $word.ActiveDocument

# Each problem row in the table has 5 cells; rows alternate between a
# row of multiplication problems and several blank spacer rows.
# Map: table row index (1-based) -> list of (old, new) pairs for its 5 cells.
$replacements = @{
    1  = @("73×16=", "59×44=", "17×61=", "77×72=", "90×29=", "55×93=", "24×83=", "25×50=", "29×83=", "68×32=")
    5  = @("37×41=", "14×62=", "57×52=", "93×56=", "76×68=", "74×98=", "46×79=", "12×83=", "52×54=", "68×73=")
    10 = @("87×63=", "90×57=", "83×60=", "82×12=", "34×87=", "12×64=", "93×41=", "58×30=", "24×56=", "92×31=")
    15 = @("29×83=", "45×89=", "38×96=", "96×96=", "89×27=", "78×54=", "63×59=", "22×37=", "55×19=", "42×44=")
    20 = @("81×98=", "70×85=", "87×77=", "17×65=", "81×19=", "88×96=", "32×93=", "20×59=", "40×26=", "49×89=")
}

$table = $d.Tables.Item(1)

foreach ($rowIndex in $replacements.Keys) {
    $pairs = $replacements[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $old = $pairs[($col - 1) * 2]
        $new = $pairs[($col - 1) * 2 + 1]
        $cell = $table.Cell($rowIndex, $col)
        $rng = $cell.Range
        $rng.MoveEnd(1, -1) | Out-Null
        if ($rng.Text -eq $old) {
            $rng.Text = $new
        }
    }
}
